# Generate Report for Handoff
#
# A fresh handoff report was generated for the four files whose status rows
# were still showing the previous run (010ff5df…, 0da89829…, 49de4360…,
# db48715b…) on both the "zh-cn" and "de-de" status sheets. The new report
# bumps each of those rows' Priority from "low" to "ht" and refreshes the
# Latest Handoff Datetime by about 30 seconds.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in @(4,5,6,7)) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-09-04 20:37:06"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in @(4,5,6,7)) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-09-04 20:37:12"
}
